$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Ligand-Receptor (Myoc-Fzd3) data rows 2-16, reflecting new TPM-based
# NATMI computation. A new "Inflammatory-Mac" cluster is introduced, and the
# former duplicate Resolving-Mac/Resolving-Mac row (old row 17) is removed.
$data = @(
    @('ECs', 'Myoc', 'Fzd3', 'ECs', 3, 1, 1.11167, 3.33501, 0.1644205457720283, 0.1644205457720283, 3, 1, 0.274713, 0.824139, 0.1055967877339779, 0.1055967877339779, 0.3053902007099999, 2.74851180639, 0.01736228147099367, 0.01736228147099367),
    @('ECs', 'Myoc', 'Fzd3', 'FAPs', 3, 1, 1.11167, 3.33501, 0.1644205457720283, 0.1644205457720283, 3, 1, 0.8886716666666666, 2.666015, 0.3415960415058637, 0.3415960415058638, 0.9879096316833332, 8.89118668515, 0.05616540757795854, 0.05616540757795856),
    @('ECs', 'Myoc', 'Fzd3', 'MuSCs', 3, 1, 1.11167, 3.33501, 0.1644205457720283, 0.1644205457720283, 3, 1, 1.438143333333333, 4.31443, 0.5528071707601584, 0.5528071707601584, 1.598740799366666, 14.3886671943, 0.09089285672307609, 0.09089285672307611),
    @('FAPs', 'Myoc', 'Fzd3', 'ECs', 3, 1, 5.135293, 15.405879, 0.7595308659577722, 0.7595308659577723, 3, 1, 0.274713, 0.824139, 0.1055967877339779, 0.1055967877339779, 1.410731745909, 12.696585713181, 0.08020401962994728, 0.08020401962994729),
    @('FAPs', 'Myoc', 'Fzd3', 'FAPs', 3, 1, 5.135293, 15.405879, 0.7595308659577722, 0.7595308659577723, 3, 1, 0.8886716666666666, 2.666015, 0.3415960415058637, 0.3415960415058638, 4.563589389131667, 41.07230450218499, 0.2594527372126957, 0.2594527372126958),
    @('FAPs', 'Myoc', 'Fzd3', 'MuSCs', 3, 1, 5.135293, 15.405879, 0.7595308659577722, 0.7595308659577723, 3, 1, 1.438143333333333, 4.31443, 0.5528071707601584, 0.5528071707601584, 7.385287392663333, 66.46758653396999, 0.4198741091151291, 0.4198741091151292),
    @('Inflammatory-Mac', 'Myoc', 'Fzd3', 'ECs', 1, 0.3333333333333333, 0.08681433333333333, 0.260443, 0.01284019544244376, 0.01284019544244376, 3, 1, 0.274713, 0.824139, 0.1055967877339779, 0.1055967877339779, 0.023849025953, 0.214641233577, 0.001355883392598524, 0.001355883392598524),
    @('Inflammatory-Mac', 'Myoc', 'Fzd3', 'FAPs', 1, 0.3333333333333333, 0.08681433333333333, 0.260443, 0.01284019544244376, 0.01284019544244376, 3, 1, 0.8886716666666666, 2.666015, 0.3415960415058637, 0.3415960415058638, 0.07714943829388889, 0.6943449446449999, 0.004386159935300421, 0.004386159935300422),
    @('Inflammatory-Mac', 'Myoc', 'Fzd3', 'MuSCs', 1, 0.3333333333333333, 0.08681433333333333, 0.260443, 0.01284019544244376, 0.01284019544244376, 3, 1, 1.438143333333333, 4.31443, 0.5528071707601584, 0.5528071707601584, 0.1248514547211111, 1.12366309249, 0.007098152114544815, 0.007098152114544816),
    @('MuSCs', 'Myoc', 'Fzd3', 'ECs', 2, 0.6666666666666666, 0.34273, 1.02819, 0.05069117062837646, 0.05069117062837647, 3, 1, 0.274713, 0.824139, 0.1055967877339779, 0.1055967877339779, 0.09415238648999999, 0.8473714784099999, 0.005352824784831524, 0.005352824784831524),
    @('MuSCs', 'Myoc', 'Fzd3', 'FAPs', 2, 0.6666666666666666, 0.34273, 1.02819, 0.05069117062837646, 0.05069117062837647, 3, 1, 0.8886716666666666, 2.666015, 0.3415960415058637, 0.3415960415058638, 0.3045744403166666, 2.741169962849999, 0.01731590322595171, 0.01731590322595171),
    @('MuSCs', 'Myoc', 'Fzd3', 'MuSCs', 2, 0.6666666666666666, 0.34273, 1.02819, 0.05069117062837646, 0.05069117062837647, 3, 1, 1.438143333333333, 4.31443, 0.5528071707601584, 0.5528071707601584, 0.4928948646333333, 4.436053781699999, 0.02802244261759323, 0.02802244261759324),
    @('Resolving-Mac', 'Myoc', 'Fzd3', 'ECs', 2, 0.6666666666666666, 0.08463066666666667, 0.253892, 0.01251722219937926, 0.01251722219937926, 3, 1, 0.274713, 0.824139, 0.1055967877339779, 0.1055967877339779, 0.023249144332, 0.209242298988, 0.001321778455606887, 0.001321778455606887),
    @('Resolving-Mac', 'Myoc', 'Fzd3', 'FAPs', 2, 0.6666666666666666, 0.08463066666666667, 0.253892, 0.01251722219937926, 0.01251722219937926, 3, 1, 0.8886716666666666, 2.666015, 0.3415960415058637, 0.3415960415058638, 0.07520887559777778, 0.67687988038, 0.004275833553957276, 0.004275833553957276),
    @('Resolving-Mac', 'Myoc', 'Fzd3', 'MuSCs', 2, 0.6666666666666666, 0.08463066666666667, 0.253892, 0.01251722219937926, 0.01251722219937926, 3, 1, 1.438143333333333, 4.31443, 0.5528071707601584, 0.5528071707601584, 0.1217110290622222, 1.09539926156, 0.006919610189815095, 0.006919610189815096),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}

# The old sheet had 16 data rows (17 rows incl. header); the new data only has
# 15 data rows, so delete the now-stale last row (old row 17).
$ws.Rows.Item(17).Delete()
